# Auto-generated edit script for API-List.xlsx "Updated API Lists in accordance with V3"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restructure rows (insert 2 new rows to grow the table from 12 -> 14 rows) ---
# one extra row in the "Admin" block (new "CRUD Employee" row) ...
$ws.Rows.Item(6).Insert()
# ... and one extra row in the "Core" block (Deliver Shipment splits into two API rows)
$ws.Rows.Item(13).Insert()

# --- Step 2: helper to apply one of the sheet's 5 cell-style archetypes ---
function Set-CellStyle {
    param($rng, [string]$kind)
    $rng.Borders.LineStyle = 1
    if ($kind -eq "bold" -or $kind -eq "boldwrap") {
        $rng.Font.Bold = $true
    } else {
        $rng.Font.Bold = $false
    }
    if ($kind -eq "boldwrap" -or $kind -eq "wrap") {
        $rng.WrapText = $true
    } else {
        $rng.WrapText = $false
    }
    if ($kind -eq "center") {
        $rng.HorizontalAlignment = -4108
        $rng.VerticalAlignment = -4108
    } else {
        $rng.HorizontalAlignment = 1
        $rng.VerticalAlignment = -4107
    }
}

# --- Step 3: write every cell value + style, row by row, for the full A1:G14 grid ---
# Row 1
$ws.Range("A1").Value = ""
Set-CellStyle $ws.Range("A1") "bold"
$ws.Range("B1").Value = 'API'
Set-CellStyle $ws.Range("B1") "bold"
$ws.Range("C1").Value = 'Requirements'
Set-CellStyle $ws.Range("C1") "boldwrap"
$ws.Range("D1").Value = 'Input '
Set-CellStyle $ws.Range("D1") "boldwrap"
$ws.Range("E1").Value = 'Output'
Set-CellStyle $ws.Range("E1") "boldwrap"
$ws.Range("F1").Value = 'Status'
Set-CellStyle $ws.Range("F1") "boldwrap"
$ws.Range("G1").Value = 'Remarks'
Set-CellStyle $ws.Range("G1") "boldwrap"

# Row 2
$ws.Range("A2").Value = 'Admin'
Set-CellStyle $ws.Range("A2") "center"
$ws.Range("B2").Value = 'CRUD Branch'
Set-CellStyle $ws.Range("B2") "plain"
$ws.Range("C2").Value = 'Create, Read, Update, Delete on Branch; 1 Stored Proc for Create and Update; While insertign into branch, make sure to insert branch''s pincode to service_pincode table as well'
Set-CellStyle $ws.Range("C2") "wrap"
$ws.Range("D2").Value = ""
Set-CellStyle $ws.Range("D2") "wrap"
$ws.Range("E2").Value = ""
Set-CellStyle $ws.Range("E2") "wrap"
$ws.Range("F2").Value = ""
Set-CellStyle $ws.Range("F2") "plain"
$ws.Range("G2").Value = 'Delete can be skipped'
Set-CellStyle $ws.Range("G2") "plain"

# Row 3
$ws.Range("A3").Value = ""
Set-CellStyle $ws.Range("A3") "center"
$ws.Range("B3").Value = 'CRUD Agent'
Set-CellStyle $ws.Range("B3") "plain"
$ws.Range("C3").Value = 'Create, Read, Update, Delete on Agent; 1 Stored Proc for Create and Update'
Set-CellStyle $ws.Range("C3") "wrap"
$ws.Range("D3").Value = ""
Set-CellStyle $ws.Range("D3") "wrap"
$ws.Range("E3").Value = ""
Set-CellStyle $ws.Range("E3") "wrap"
$ws.Range("F3").Value = ""
Set-CellStyle $ws.Range("F3") "plain"
$ws.Range("G3").Value = 'Delete can be skipped'
Set-CellStyle $ws.Range("G3") "plain"

# Row 4
$ws.Range("A4").Value = ""
Set-CellStyle $ws.Range("A4") "center"
$ws.Range("B4").Value = 'CRUD Employee'
Set-CellStyle $ws.Range("B4") "plain"
$ws.Range("C4").Value = 'Create, Read, Update, Delete on Employee; 1 Stored Proc for Create and Update; Create a sequence to generate employee_id of 6 digits'
Set-CellStyle $ws.Range("C4") "wrap"
$ws.Range("D4").Value = ""
Set-CellStyle $ws.Range("D4") "wrap"
$ws.Range("E4").Value = ""
Set-CellStyle $ws.Range("E4") "wrap"
$ws.Range("F4").Value = ""
Set-CellStyle $ws.Range("F4") "plain"
$ws.Range("G4").Value = 'Delete can be skipped'
Set-CellStyle $ws.Range("G4") "plain"

# Row 5
$ws.Range("A5").Value = ""
Set-CellStyle $ws.Range("A5") "center"
$ws.Range("B5").Value = 'CRUD Customer'
Set-CellStyle $ws.Range("B5") "plain"
$ws.Range("C5").Value = 'Create, Read, Update, Delete on Customer; 1 Stored Proc for Create and Update'
Set-CellStyle $ws.Range("C5") "wrap"
$ws.Range("D5").Value = ""
Set-CellStyle $ws.Range("D5") "wrap"
$ws.Range("E5").Value = ""
Set-CellStyle $ws.Range("E5") "wrap"
$ws.Range("F5").Value = ""
Set-CellStyle $ws.Range("F5") "plain"
$ws.Range("G5").Value = 'Delete can be skipped'
Set-CellStyle $ws.Range("G5") "plain"

# Row 6
$ws.Range("A6").Value = ""
Set-CellStyle $ws.Range("A6") "center"
$ws.Range("B6").Value = 'CRUD Serviceable Pincodes'
Set-CellStyle $ws.Range("B6") "plain"
$ws.Range("C6").Value = 'Create, Read, Update, Delete on Serviceable Pincodes'
Set-CellStyle $ws.Range("C6") "wrap"
$ws.Range("D6").Value = ""
Set-CellStyle $ws.Range("D6") "wrap"
$ws.Range("E6").Value = ""
Set-CellStyle $ws.Range("E6") "wrap"
$ws.Range("F6").Value = ""
Set-CellStyle $ws.Range("F6") "plain"
$ws.Range("G6").Value = 'Can be skipped. Insert scripts are sufficient'
Set-CellStyle $ws.Range("G6") "plain"

# Row 7
$ws.Range("A7").Value = ""
Set-CellStyle $ws.Range("A7") "plain"
$ws.Range("B7").Value = ""
Set-CellStyle $ws.Range("B7") "plain"
$ws.Range("C7").Value = ""
Set-CellStyle $ws.Range("C7") "wrap"
$ws.Range("D7").Value = ""
Set-CellStyle $ws.Range("D7") "wrap"
$ws.Range("E7").Value = ""
Set-CellStyle $ws.Range("E7") "wrap"
$ws.Range("F7").Value = ""
Set-CellStyle $ws.Range("F7") "plain"
$ws.Range("G7").Value = ""
Set-CellStyle $ws.Range("G7") "plain"

# Row 8
$ws.Range("A8").Value = 'Core'
Set-CellStyle $ws.Range("A8") "center"
$ws.Range("B8").Value = 'Get Quote'
Set-CellStyle $ws.Range("B8") "plain"
$ws.Range("C8").Value = 'Customer can check whether the service is available between two given routes`nOptional - also get estimate cost'
Set-CellStyle $ws.Range("C8") "wrap"
$ws.Range("D8").Value = 'source pin, dest pin, shipment dimensions (optional)'
Set-CellStyle $ws.Range("D8") "wrap"
$ws.Range("E8").Value = 'boolean, cost'
Set-CellStyle $ws.Range("E8") "wrap"
$ws.Range("F8").Value = 'Complete'
Set-CellStyle $ws.Range("F8") "wrap"
$ws.Range("G8").Value = ""
Set-CellStyle $ws.Range("G8") "plain"

# Row 9
$ws.Range("A9").Value = ""
Set-CellStyle $ws.Range("A9") "center"
$ws.Range("B9").Value = 'Find Route'
Set-CellStyle $ws.Range("B9") "plain"
$ws.Range("C9").Value = 'Agent finds route and price details based on shipment''s details'
Set-CellStyle $ws.Range("C9") "wrap"
$ws.Range("D9").Value = 'source pin, dest pin, shipment dimensions, other shipment details'
Set-CellStyle $ws.Range("D9") "wrap"
$ws.Range("E9").Value = 'route and cost'
Set-CellStyle $ws.Range("E9") "wrap"
$ws.Range("F9").Value = 'Complete'
Set-CellStyle $ws.Range("F9") "wrap"
$ws.Range("G9").Value = ""
Set-CellStyle $ws.Range("G9") "plain"

# Row 10
$ws.Range("A10").Value = ""
Set-CellStyle $ws.Range("A10") "center"
$ws.Range("B10").Value = 'Book shipment'
Set-CellStyle $ws.Range("B10") "plain"
$ws.Range("C10").Value = 'Once the payment is received (offline), Agent receives the shipment and performs below tasks:`n1. Checks if customer is registered looking at his contact_num. If no, then registers the customer. `n2. A record is inserted into shipment:`n    Update its attributes such that cusotmer_id=regd customer id (source branch and address details can be  determined using customer''s address details), dest_branch=branch of dest pincode, next_branch = source_branch; status = booked`n3. Insert a record into shipment_tracker shipment_id; agent_id, current_branch = null; next_branch=source_branch; creation_datetime, status=booked'
Set-CellStyle $ws.Range("C10") "wrap"
$ws.Range("D10").Value = 'shipment details + cost'
Set-CellStyle $ws.Range("D10") "wrap"
$ws.Range("E10").Value = 'shipment_id and status'
Set-CellStyle $ws.Range("E10") "wrap"
$ws.Range("F10").Value = ""
Set-CellStyle $ws.Range("F10") "plain"
$ws.Range("G10").Value = ""
Set-CellStyle $ws.Range("G10") "plain"

# Row 11
$ws.Range("A11").Value = ""
Set-CellStyle $ws.Range("A11") "center"
$ws.Range("B11").Value = 'Receive Shipment at branch'
Set-CellStyle $ws.Range("B11") "plain"
$ws.Range("C11").Value = 'When a shipment is arrived at any branch, employee at the branch will perform below tasks:`n1. Update shipment table with status = RECEIVED_AT_DEST_BRANCH if current_branch = dest_branch else status = IN_TRANSIT`n2. Insert a row into shipment_tracker table by referring to route table:`n shipment_id; employee_id, current_branch = receiving branch; next_branch=determine from route detail; creation_datetime, status=RECEIVED_AT_DEST_BRANCH or IN_TRANSIT'
Set-CellStyle $ws.Range("C11") "wrap"
$ws.Range("D11").Value = 'receiving branch code, shipment_id'
Set-CellStyle $ws.Range("D11") "wrap"
$ws.Range("E11").Value = 'shipment_id and status'
Set-CellStyle $ws.Range("E11") "wrap"
$ws.Range("F11").Value = ""
Set-CellStyle $ws.Range("F11") "plain"
$ws.Range("G11").Value = ""
Set-CellStyle $ws.Range("G11") "plain"

# Row 12
$ws.Range("A12").Value = ""
Set-CellStyle $ws.Range("A12") "center"
$ws.Range("B12").Value = 'Attempt Delivery'
Set-CellStyle $ws.Range("B12") "plain"
$ws.Range("C12").Value = 'When a shipment is marked as RECEIVED_AT_DEST_BRANCH, an agent from the destination branch attempts the delivery and performs below task:`n1. Updates shipment table such that status = OUT_FOR_DELIVERY`n2. Insert a record into shipment_tracker; shipment_id, agent_id, creation_datetime,status=OUT_FOR_DELIVERY`n'
Set-CellStyle $ws.Range("C12") "wrap"
$ws.Range("D12").Value = 'shipment_id'
Set-CellStyle $ws.Range("D12") "wrap"
$ws.Range("E12").Value = 'shipment_id and status'
Set-CellStyle $ws.Range("E12") "wrap"
$ws.Range("F12").Value = ""
Set-CellStyle $ws.Range("F12") "plain"
$ws.Range("G12").Value = ""
Set-CellStyle $ws.Range("G12") "plain"

# Row 13
$ws.Range("A13").Value = ""
Set-CellStyle $ws.Range("A13") "center"
$ws.Range("B13").Value = 'Update Delivery Status'
Set-CellStyle $ws.Range("B13") "plain"
$ws.Range("C13").Value = 'When a shipment is marked as OUT_FOR_DELIVERY, the same agent from the destination branch updates the status of the delivery and performs below task:`n1. Updates shipment table such that status = DELIVERED or UNDELIVERED and adds status_remarks (optional)`n2. Insert a record into shipment_tracker; shipment_id, agent_id, creation_datetime,status=DELIVERED or UNDELIVERED`n'
Set-CellStyle $ws.Range("C13") "wrap"
$ws.Range("D13").Value = 'shipment_id'
Set-CellStyle $ws.Range("D13") "wrap"
$ws.Range("E13").Value = 'shipment_id and status'
Set-CellStyle $ws.Range("E13") "wrap"
$ws.Range("F13").Value = ""
Set-CellStyle $ws.Range("F13") "plain"
$ws.Range("G13").Value = ""
Set-CellStyle $ws.Range("G13") "plain"

# Row 14
$ws.Range("A14").Value = ""
Set-CellStyle $ws.Range("A14") "center"
$ws.Range("B14").Value = 'Track Shipment'
Set-CellStyle $ws.Range("B14") "plain"
$ws.Range("C14").Value = 'At any point in time, customer should be able to track the full history of the shipment. This API would query shipment_tracker table and produce required output'
Set-CellStyle $ws.Range("C14") "wrap"
$ws.Range("D14").Value = 'shipment_id'
Set-CellStyle $ws.Range("D14") "wrap"
$ws.Range("E14").Value = 'shipment and shipment_tracker'
Set-CellStyle $ws.Range("E14") "wrap"
$ws.Range("F14").Value = ""
Set-CellStyle $ws.Range("F14") "plain"
$ws.Range("G14").Value = ""
Set-CellStyle $ws.Range("G14") "plain"

# --- Step 4: row heights (match target autofit heights for wrapped multi-line rows) ---
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 115.2
$ws.Rows.Item(11).RowHeight = 86.4
$ws.Rows.Item(12).RowHeight = 86.4
$ws.Rows.Item(13).RowHeight = 100.8
$ws.Rows.Item(14).RowHeight = 28.8

# --- Step 5: merged cells for the two "API group" label columns ---
$ws.Range("A2:A5").UnMerge()
$ws.Range("A8:A14").UnMerge()
$ws.Range("A2:A6").Merge()
$ws.Range("A8:A14").Merge()

# --- Step 6: column widths for the now-repurposed F (Status) / G (Remarks) columns ---
$ws.Columns.Item(6).ColumnWidth = 8
$ws.Columns.Item(7).ColumnWidth = 34.666666666666664

# --- Step 7: selection / active cell to match the saved view state ---
$ws.Range("E14").Select()
